$wb = $excel.ActiveWorkbook

# Sheet 1: OFF
$ws1 = $wb.Worksheets.Item("OFF")
$ws1.Range("B2").Value = 267
$ws1.Range("C2").Value = 199
$ws1.Range("D2").Value = 61
$ws1.Range("E2").Value = 32

# Sheet 2: DEF
$ws2 = $wb.Worksheets.Item("DEF")
$ws2.Range("B2").Value = 313
$ws2.Range("C2").Value = 213
$ws2.Range("D2").Value = 74
$ws2.Range("E2").Value = 23
